$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4: add "1 - x" formulas for the complementary probability cells ---
$ws.Range("D4").Formula = "=1-C4"
$ws.Range("H4").Formula = "=1-G4"
$ws.Range("L4").Formula = "=1-K4"

# --- Row 5: same pattern ---
$ws.Range("D5").Formula = "=1-C5"
$ws.Range("H5").Formula = "=1-G5"
$ws.Range("L5").Formula = "=1-K5"

# --- Row 12: swap which column holds the literal vs. the formula ---
$ws.Range("D12").Value = 0.13
$ws.Range("E12").Formula = "=1-D12"

# --- Rows 13-15, 17-20, 22-25: replace the literal probabilities in column D ---
$ws.Range("D13").Value = 0.58
$ws.Range("D14").Value = 0.15
$ws.Range("D15").Value = 0.6

$ws.Range("D17").Value = 0.15
$ws.Range("D18").Value = 0.55
$ws.Range("D19").Value = 0.15
$ws.Range("D20").Value = 0.63

$ws.Range("D22").Value = 0.15
$ws.Range("D23").Value = 0.66
$ws.Range("D24").Value = 0.15
$ws.Range("D25").Value = 0.56

# --- Update the saved cell selection to match the author's final position ---
$ws.Range("G15").Select()
